# Applies the RQ1/RQ2 wording tweaks on the "STGraph - Conceptualization"
# slide and swaps the ASCII "<->" arrows for the Unicode "⟺" glyph on the
# "Graph+TimeSeries Hybrid data model" slide, while preserving the existing
# run/paragraph structure (bold "RQ1"/"RQ2" labels, bullet levels, etc.).

function Replace-InTextRange {
    param($TextRange, $OldText, $NewText)

    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -ge 0) {
        $sub = $TextRange.Characters($idx + 1, $OldText.Length)
        $sub.Text = $NewText
    }
}

$p = $ppt.ActivePresentation

# --- Slide 3: "STGraph - Conceptualization" ------------------------------
$s3 = $p.Slides.Item(3)
$rqShape = $s3.Shapes.Item(2)
$rqRange = $rqShape.TextFrame.TextRange

Replace-InTextRange $rqRange " - Can we separate data by temporal granularity to support hybrid modeling across time-series and temporal graph systems?" " - How can we separate, within a specific domain, temporal graph data from time-series data?"

Replace-InTextRange $rqRange " - If so, can we embed two different data-layout into the same conceptual storage system and provide hybrid capabilities?" " - If so, can we embed both data layouts into the same conceptual storage system and provide hybrid capabilities?"

# --- Slide 4: "Graph+TimeSeries Hybrid data model" ------------------------
$s4 = $p.Slides.Item(4)
$edgeShape = $s4.Shapes.Item(3)
$edgeRange = $edgeShape.TextFrame.TextRange

Replace-InTextRange $edgeRange "Graph edge  <–> " "Graph edge  ⟺ "

Replace-InTextRange $edgeRange "Virtual edge  <–> " "Virtual edge  ⟺ "
